# Adds 13 new abituryent (applicant) rows for the 06/29-07/01 intake to the
# "Лист1" sheet (rows 128-140), extending the sheet dimension accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,@('FARMONOVA ZARINA','Psixologiya','O''zbek tili','Kunduzgi','AD5184165','60109076960015','Toshkent shahri','Yashnaobod tumani','998911623523','+998918102334','2025-06-29')
    ,@('Eshmuratov Jaloliddin Tohir o''g''li','Yurisprudensiya','O''zbek tili','Kunduzgi','AD0610768','50909045470038','Jizzax viloyati','Zomin tumani','998949474230','+998957600076','2025-06-29')
    ,@('Halilov Abdurauf','Yurisprudensiya','Rus tili','Kunduzgi','AD3879501','32703922200035','Toshkent viloyati','Angren shahri','998937355565','+998937355565','2025-06-29')
    ,@('Oktamova ruhshona','Yurisprudensiya','O''zbek tili','Kunduzgi','AD6121348','60810067170055','Xorazm viloyati','Shovot tumani','998939140812','+998943106633','2025-06-29')
    ,@('Bu qaysi harf','Yurisprudensiya','O''zbek tili','Kunduzgi','AD2755823','32904931690024','Jizzax viloyati','Yangiobod tumani','998955779302','+998955779302','2025-06-30')
    ,@('Kim Yuriy','Mehnat muhofazasi va texnika xavfsizligi','Rus tili','Kunduzgi','AD3569381','32008781930014','Toshkent shahri','Mirzo Ulugʻbek tumani','+998977337388','+998977337388','2025-06-30')
    ,@('Boratov Vohidjon Erkinboy o''g''li','Yurisprudensiya','O''zbek tili','Kunduzgi','AD4095733','30302924280027','Fargona viloyati','Dangʻara tumani','+998905868787','+998500117747','2025-06-30')
    ,@('Rajabova Iroda Obid qizi','Yurisprudensiya','O''zbek tili','Kunduzgi','AD3282682','62608076060115','Samarqand viloyati','Kattaqoʻrgʻon tumani','998940149199','+998947535885','2025-06-30')
    ,@('Nishonva Dildora Ziyodulla qizi','Yurisprudensiya','O''zbek tili','Kunduzgi','AC2650600','41606956420012','Sirdaryo viloyati','Sardoba tumani','998979061122','+998979061122','2025-06-30')
    ,@('Toshpulatov Jamshid Xodjiakbar o''g''li','Yurisprudensiya','Rus tili','Kunduzgi','AD3351907','31104920270760','Toshkent shahri','Yashnaobod tumani','998935477271','+998935477271','2025-06-30')
    ,@('Maqsudov Nurmuhammad shavkatjon ogli','Yurisprudensiya','O''zbek tili','Kunduzgi','AD6260332','51908065330644','Buxoro viloyati','Peshku tumani','998953540772','+998996846250','2025-06-30')
    ,@('Zakiryayeva Shirin Nurillayevna','Yurisprudensiya','Rus tili','Kunduzgi','AD1874714','41310826180015','Samarqand viloyati','Samarqand shahri','998979244447','+998979244447','2025-07-01')
    ,@('Ochilova Ruxshona Abdurashidovna','Yurisprudensiya','O''zbek tili','Kunduzgi','AD3963854','61403076150025','Samarqand viloyati','Toyloq tumani','998507147607','+998507147607','2025-07-01')
)

$startRow = 128
$endRow = $startRow + $newRows.Length - 1

# JSHIR (F), phone/telegram numbers (I, J) and the yyyy-mm-dd enrollment
# date (K) all look numeric/date-like to Excel, so pre-format those
# columns as text before writing the values - otherwise they would be
# silently coerced into numbers or date serials, losing leading "+"
# signs and the original "yyyy-mm-dd" text representation.
$ws.Range("F$startRow`:F$endRow").NumberFormat = "@"
$ws.Range("I$startRow`:I$endRow").NumberFormat = "@"
$ws.Range("J$startRow`:J$endRow").NumberFormat = "@"
$ws.Range("K$startRow`:K$endRow").NumberFormat = "@"

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowNum = $startRow + $r
    $rowValues = $newRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}
